$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0, "falling", 0.5757570266723633, -0.6172752380371094, -0.5019410252571106, -0.005192354787141, 0.0755945742130279, 0.0082466807216405),
    @(100, "falling", 0.4425497055053711, -0.6991405487060547, -0.699306845664978, 0.00167987938039, 0.0209221355617046, -0.0082466807216405),
    @(200, "falling", 0.3754444122314453, -0.6968369483947754, -0.6064528226852417, 0.0074830991216003, -0.0106901414692401, 0.0113010071218013),
    @(300, "falling", 0.4661340713500976, -0.7378168106079102, -0.8263083696365356, -0.0024434609804302, 0.0310014113783836, 0.0187841057777404),
    @(400, "falling", 0.3465394973754883, -0.7457756996154785, -0.6906525492668152, -0.0039706239476799, 0.0114537235349416, -0.0603229440748691),
    @(500, "falling", 0.1021490097045898, -0.7542791366577148, -0.5537225604057312, -0.08918632566928859, 0.2981022298336029, 0.0048869219608604),
    @(600, "falling", 0.2407388687133789, -0.6835846900939941, -0.2614910900592804, -0.2064724564552307, 0.7906123399734497, 0.058184914290905),
    @(700, "falling", 0.7179374694824219, -0.8143949508666992, 0.3482142686843872, -0.294895201921463, 0.9447031021118164, 0.0591012127697467),
    @(800, "falling", 0.5368337631225586, -1.205713748931885, 1.448910713195801, -0.011148290708661, -0.4109596014022827, 0.3446807265281677),
    @(900, "falling", 0.1711635589599609, -1.118285655975342, 2.145041465759277, -0.3101668357849121, -0.0274889357388019, -0.0977384373545646),
    @(1000, "falling", -0.2017207145690918, -1.242213249206543, 2.87759017944336, -0.2105957865715026, 0.9750936627388, -1.285871386528015),
    @(1100, "falling", -1.174193859100342, -2.213449716567993, 2.043629169464112, 0.1690569519996643, -0.0616973899304866, -0.9847147464752196),
    @(1200, "falling", -2.267188549041748, -3.465161085128784, 2.441758394241333, -1.828319668769836, 1.374752283096314, 2.351678371429444),
    @(1300, "falling", -3.827628612518311, -2.444296360015869, 1.300557613372803, 1.769065737724304, -1.171639561653137, -0.7515169382095337),
    @(1400, "falling", -4.013436794281006, -1.781657457351685, 1.157623410224915, 0.3782783150672912, -0.1327104717493057, 0.1701259762048721),
    @(1500, "falling", -6.23094367980957, -1.233789682388306, -1.820952653884888, -0.3932445049285888, -0.6039929986000061, 1.466076612472534),
    @(1600, "falling", 23.2268295288086, -21.41561698913575, -21.18720245361328, -0.2658790946006775, -0.8356636762619019, -1.123686671257019),
    @(1700, "falling", -2.458966255187988, -0.9571634531021118, -1.108598947525024, -0.1258382350206375, -0.4847215712070465, -0.2483167201280594),
    @(1800, "falling", 1.693714141845703, -1.90864372253418, 0.4882011413574219, -0.1476766765117645, 0.0491746515035629, 0.3875939846038818),
    @(1900, "falling", -0.7738790512084961, -0.2492363452911377, -2.26715350151062, 0.001527163083665, 0.3072652220726013, 0.0675006061792373),
    @(2000, "falling", 5.673647880554199, 0.3179191350936889, 1.958111882209778, -0.0311541277915239, 0.040775254368782, -0.4699080884456634),
    @(2100, "falling", 0.6818609237670898, -1.483997821807861, 1.116373777389526, 0.007941247895359899, -0.0103847095742821, 0.0820086598396301),
    @(2200, "falling", -1.090703964233398, -0.7039613127708435, 0.7109836339950562, -0.012980886735022, 0.1769981980323791, 0.0358883328735828),
    @(2300, "falling", 1.377178192138672, 0.5982787609100342, -0.2494584619998932, 0.0025961773935705, -0.0377209298312664, -0.1343903541564941),
    @(2400, "falling", -1.309194564819336, -1.633804559707642, -0.4874744415283203, -0.0339030213654041, -0.0348193198442459, 0.0242818929255008),
    @(2500, "falling", 0.3532924652099609, -1.432539701461792, -0.3063055276870727, -0.0064140851609408, 0.0251981914043426, 0.0591012127697467),
    @(2600, "falling", 0.4603271484375, -0.7004889249801636, -0.5307044386863708, 0.0099265603348612, -0.0109955742955207, -0.0491746515035629),
    @(2700, "falling", -0.4794178009033203, -1.04423999786377, -0.2169336676597595, -0.0039706239476799, 0.0164933614432811, -0.0076358155347406),
    @(2800, "falling", 0.1147146224975586, -0.7852307558059692, -0.0405309796333313, -0.0065668015740811, -0.0163406450301408, 0.0088575463742017),
    @(2900, "falling", -0.0877876281738281, -0.7322115302085876, -0.3498360514640808, -0.009315694682300001, -0.0369573459029197, -0.0157297793775796),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
    $ws.Range("H$r").Value = $row[7]
}
